$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp text in F1
$ws.Range("F1").Value = "Last status check on: 21.02.2022 00:30"

# Update D10 from the inline text "+0.9" to the numeric value 0.9
$ws.Range("D10").Value = 0.9

# Update E10 from the inline text "2022-02-21 00:21:04" to the numeric
# Excel date serial value, formatted the same way as the other date cells
# in column E (e.g. E2, E3, ...).
$ws.Range("E10").Value = 44613.01462962963
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
